$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set F3:F82 to 1 (adds themes/cultural-context vocabulary marker column)
$ws.Range("F3:F82").Value = 1

# Update the selection to match the new active range
$ws.Range("F6:F82").Select()
